$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so that
# numeric-looking strings (e.g. "0.600", "41.60", "51.738.24")
# are preserved exactly as text instead of being coerced to numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "51.738.24"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.820.33"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "351.98"
$ws.Range("E5").Value = "  +5.84%  "
$ws.Range("D6").Value = "113.57"
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  +2.88%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("D10").Value = "41.60"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").Value = "0.0850"
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "19.98"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.131"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "7.72"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "3.262.67"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").Value = "2.802.64"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "0.893"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "51.587.46"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  +7.36%  "
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").Value = "13.49"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "0.0₃0994"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "269.93"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "2.77"
$ws.Range("E25").Value = "  +4.30%  "
$ws.Range("D26").Value = "26.67"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("D31").Value = "50.61"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "33.73"
$ws.Range("E32").Value = "  -3.90%  "
$ws.Range("D33").Value = "0.0450"
$ws.Range("E33").Value = "  +27.57%  "
$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  +4.33%  "
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D40").Value = "18.10"
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").Value = "23.73"
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").Value = "2.56"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "125.42"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "2.081.56"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "0.939"
$ws.Range("E50").Value = "  +7.18%  "
$ws.Range("D51").Value = "60.67"
$ws.Range("E51").Value = "  +1.10%  "

# Restore the original (default) cell style now that the text
# values are safely stored, so formatting matches the source.
$priceRange.Style = "Normal"

